$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to be stored as text so that numeric-looking
# values (e.g. "390.33") are not silently converted into Excel numbers and
# lose their original formatting / exact text representation.
$ws.Range("D2:D51").NumberFormat = "@"

$changes = @{
    "D2" = '51.821.87'
    "E2" = '  +0.51%  '
    "D3" = '3.102.34'
    "E3" = '  +3.82%  '
    "E4" = '  +0.07%  '
    "D5" = '390.33'
    "E5" = '  +2.28%  '
    "D6" = '103.70'
    "E6" = '  -0.39%  '
    "D7" = '0.545'
    "E7" = '  -0.12%  '
    "E8" = '  +0.08%  '
    "D9" = '0.592'
    "E9" = '  -0.59%  '
    "D10" = '37.27'
    "E10" = '  +1.43%  '
    "E11" = '  -0.11%  '
    "D12" = '0.0863'
    "E12" = '  +0.39%  '
    "D13" = '3.600.12'
    "E13" = '  +3.98%  '
    "D14" = '18.71'
    "E14" = '  +1.20%  '
    "D15" = '7.87'
    "E15" = '  +0.53%  '
    "D16" = '3.112.47'
    "E16" = '  +4.39%  '
    "E17" = '  -0.81%  '
    "E18" = '  -4.03%  '
    "D19" = '51.926.46'
    "E19" = '  +0.66%  '
    "D20" = '3.21'
    "E20" = '  +4.26%  '
    "E21" = '  -0.55%  '
    "E22" = '  +0.46%  '
    "D23" = '70.09'
    "E23" = '  -0.52%  '
    "D24" = '268.13'
    "E24" = '  +0.32%  '
    "E25" = '  -2.31%  '
    "E26" = '  +2.75%  '
    "E27" = '  +4.18%  '
    "E28" = '  +0.04%  '
    "B29" = 'Dai'
    "C29" = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    "D29" = '1.00'
    "E29" = '  +0.11%  '
    "B30" = 'RenderToken'
    "C30" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    "D30" = '7.18'
    "E30" = '  -0.31%  '
    "E31" = '  -0.51%  '
    "D32" = '10.36'
    "E32" = '  -0.20%  '
    "D33" = '35.57'
    "E33" = '  +2.96%  '
    "E34" = '  +3.47%  '
    "B35" = 'VeChain'
    "C35" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    "D35" = '0.0451'
    "E35" = '  +1.25%  '
    "B36" = 'OKB'
    "C36" = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    "D36" = '50.16'
    "E36" = '  -2.38%  '
    "E37" = '  -0.14%  '
    "D38" = '3.41'
    "E38" = '  +3.48%  '
    "E39" = '  +6.99%  '
    "E40" = '  +2.66%  '
    "E41" = '  +0.97%  '
    "D42" = '16.91'
    "E42" = '  -0.11%  '
    "D43" = '128.95'
    "E43" = '  -0.11%  '
    "E44" = '  -0.07%  '
    "E45" = '  -3.53%  '
    "E46" = '  +3.98%  '
    "E47" = '  +6.28%  '
    "E48" = '  +2.11%  '
    "D49" = '2.050.40'
    "E49" = '  +0.94%  '
    "D50" = '3.414.54'
    "E50" = '  +3.91%  '
    "D51" = '0.0321'
    "E51" = '  -2.61%  '
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}

# Restore the default "Normal" style on column D so no stray number-format
# style is left behind on cells that didn't originally have one.
$ws.Range("D2:D51").Style = "Normal"
